$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 36.2
$ws.Range("I8").Value = 36.2
$ws.Range("K8").Value = 108.6
$ws.Range("M8").Value = 30.39999999999999
$ws.Range("H97").Value = 1011.1111
$ws.Range("J97").Value = 1011.1111
$ws.Range("L97").Value = 3033.3333
$ws.Range("N97").Value = -4025.3333
$ws.Range("H100").Value = 474.5
$ws.Range("I100").Value = 299.33334
$ws.Range("K100").Value = 299.33334
$ws.Range("M100").Value = 241.66666
$ws.Range("H113").Value = 1600
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 1600
$ws.Range("K113").Value = 0
$ws.Range("L113").ClearContents()
$ws.Range("M113").Value = 1600
$ws.Range("N113").Value = -8108
$ws.Range("H132").Value = 3074.5
$ws.Range("I132").Value = 3074.5
$ws.Range("K132").Value = 9223.5
$ws.Range("M132").Value = -6693.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("M2").ClearContents()
$ws.Range("H74").Value = 1143.8
$ws.Range("I74").Value = 1143.8
$ws.Range("K74").Value = 1143.8
$ws.Range("M74").Value = -269.8
$ws.Range("H77").Value = 1143.8
$ws.Range("I77").Value = 1143.8
$ws.Range("K77").Value = 5719
$ws.Range("M77").Value = -1351
$ws.Range("H97").Value = 2188.7778
$ws.Range("I97").Value = 1683.1666
$ws.Range("K97").Value = 1683.1666
$ws.Range("M97").Value = -1187.1666
$ws.Range("H102").Value = 2222
$ws.Range("I102").Value = 2222
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 2222
$ws.Range("L102").Value = 0
$ws.Range("M102").ClearContents()
$ws.Range("N102").Value = -600
$ws.Range("H110").Value = 999
$ws.Range("I110").Value = 999
$ws.Range("K110").Value = 999
$ws.Range("M110").Value = 1046
$ws.Range("H116").Value = 0
$ws.Range("I116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("M116").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("M3").ClearContents()
$ws.Range("H107").Value = 1764.6471
$ws.Range("I107").Value = 1716.0714
$ws.Range("K107").Value = 1716.0714
$ws.Range("M107").Value = 203.9286
$ws.Range("H112").Value = 97999.664
$ws.Range("J112").Value = 97999.664
$ws.Range("L112").Value = 97999.664
$ws.Range("N112").Value = -100953.664
$ws.Range("H134").Value = 3747.5
$ws.Range("I134").Value = 3747.5
$ws.Range("K134").Value = 11242.5
$ws.Range("M134").Value = -8707.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1990.6666
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 1990.6666
$ws.Range("K16").Value = 0
$ws.Range("L16").ClearContents()
$ws.Range("M16").Value = 1990.6666
$ws.Range("N16").Value = -2564.6666
$ws.Range("H22").Value = 691.125
$ws.Range("I22").Value = 691.125
$ws.Range("K22").Value = 691.125
$ws.Range("M22").Value = -341.125
$ws.Range("H31").Value = 949.625
$ws.Range("J31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("N31").ClearContents()
$ws.Range("H34").Value = 949.625
$ws.Range("J34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("N34").ClearContents()
$ws.Range("H86").Value = 8425.4
$ws.Range("I86").Value = 8980.714
$ws.Range("J86").Value = 7129.6665
$ws.Range("K86").Value = 8980.714
$ws.Range("L86").Value = 7129.6665
$ws.Range("M86").Value = -7857.714
$ws.Range("N86").Value = -9375.666499999999
$ws.Range("H89").Value = 8425.4
$ws.Range("I89").Value = 8980.714
$ws.Range("J89").Value = 7129.6665
$ws.Range("K89").Value = 44903.57
$ws.Range("L89").Value = 35648.3325
$ws.Range("M89").Value = -39287.57
$ws.Range("N89").Value = -46880.3325
$ws.Range("H107").Value = 298.33334
$ws.Range("I107").Value = 323.125
$ws.Range("K107").Value = 323.125
$ws.Range("M107").Value = 1596.875
$ws.Range("H113").Value = 1990.6666
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 1990.6666
$ws.Range("K113").Value = 0
$ws.Range("L113").ClearContents()
$ws.Range("M113").Value = 1990.6666
$ws.Range("N113").Value = -6330.6666
$ws.Range("H122").Value = 1999
$ws.Range("I122").Value = 1999
$ws.Range("K122").Value = 5997
$ws.Range("M122").Value = -3547
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()
$ws.Range("H134").Value = 1991
$ws.Range("I134").Value = 1991.25
$ws.Range("K134").Value = 5973.75
$ws.Range("M134").Value = -3438.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 627.4286
$ws.Range("I33").Value = 672
$ws.Range("K33").Value = 4032
$ws.Range("M33").Value = -3749
$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("M76").ClearContents()
$ws.Range("H79").Value = 0
$ws.Range("I79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("M79").ClearContents()
$ws.Range("H109").Value = 813.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 3636
$ws.Range("I70").Value = 3499.5
$ws.Range("K70").Value = 3499.5
$ws.Range("M70").Value = -3229.5
$ws.Range("H73").Value = 3636
$ws.Range("I73").Value = 3499.5
$ws.Range("K73").Value = 3499.5
$ws.Range("M73").Value = -2563.5
$ws.Range("H113").Value = 1300
$ws.Range("J113").Value = 1300
$ws.Range("L113").Value = 1300
$ws.Range("N113").Value = -5640
$ws.Range("H132").Value = 2641.7144
$ws.Range("I132").Value = 2248.6667
$ws.Range("K132").Value = 6746.000100000001
$ws.Range("M132").Value = -4216.000100000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2115.3845
$ws.Range("I68").Value = 1972.7273
$ws.Range("K68").Value = 1972.7273
$ws.Range("M68").Value = -1223.7273
$ws.Range("H71").Value = 2115.3845
$ws.Range("I71").Value = 1972.7273
$ws.Range("K71").Value = 9863.636500000001
$ws.Range("M71").Value = -6119.636500000001
$ws.Range("H93").Value = 2965
$ws.Range("I93").Value = 2965
$ws.Range("K93").Value = 2965
$ws.Range("M93").Value = -1717

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()
